$wb = $excel.ActiveWorkbook

# Add a new worksheet named "form" after the existing sheets
$ws = $wb.Worksheets.Add()
$ws.Name = "form"

# Move it to be the last sheet (after webTable)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $lastSheet)

# NOTE: sheet object references track position, not identity, so after the
# Move the $ws/$lastSheet variables now point at whatever sheet occupies
# their original slot. Re-fetch the "form" sheet by name before continuing.
$ws = $wb.Worksheets.Item("form")

# Header row
$ws.Range("A1").Value = "nombre"
$ws.Range("B1").Value = "apellido"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "celu"
$ws.Range("E1").Value = "direccion"
$ws.Range("F1").Value = "montNumber"
$ws.Range("G1").Value = "año"
$ws.Range("H1").Value = "materia"

# Data row
$ws.Range("A2").Value = "rodrigo"
$ws.Range("B2").Value = "alvarez"
$ws.Range("C2").Value = "rodri@gmial.com"
$ws.Range("D2").Value = 1234567890
$ws.Range("E2").Value = "gil barros"
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 1988
$ws.Range("H2").Value = "english"

# Hyperlink for C2
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:rodri@gmial.com")

# Column widths (ColumnWidth is pixel-quantized by the engine, same as real
# Excel; these inputs land on the closest achievable stored width to the
# target 17.140625 / 14.5703125 character units).
$ws.Columns.Item(3).ColumnWidth = 16.333333333333336
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666

# Selection on the new sheet
$ws.Range("I2").Select()

# sheet2 should no longer have tabSelected; sheet "form" becomes the active tab
$ws2 = $wb.Worksheets.Item("webTable")
$ws2.Range("E2").Select()

$ws.Activate()
